$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the very top; existing rows (and their styles) shift down by one.
$ws.Rows.Item(1).Insert()

# Populate the new first row with the echoed command (no special styling).
$ws.Cells.Item(1, 1).Value = "#"
$ws.Cells.Item(1, 2).Value = "xltablediff.py --newSheet=Sheet2 --key=ID test1in.xlsx test1in.xlsx --out=test1out.xlsx"

# Materialize the remaining (blank) cells in the row, matching the original
# sheet's convention of an explicit empty cell per column, without adding style.
$ws.Cells.Item(1, 3).Style = "Normal"
$ws.Cells.Item(1, 4).Style = "Normal"
$ws.Cells.Item(1, 5).Style = "Normal"
$ws.Cells.Item(1, 6).Style = "Normal"
$ws.Cells.Item(1, 7).Style = "Normal"
